# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2210"
#   "<header>_new" -> "<header>_FV2304"
# and expose the data range as a proper Excel Table ("Table1"), with the
# first row frozen (as a header row) in the worksheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (columns A:J were "<name>_old", L:U were
#        "<name>_new"; column K holds the constant "diff" header). ---------
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

$leftCols  = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$rightCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $leftAddr = $leftCols[$i] + "1"
    $ws.Range($leftAddr).Value = $baseNames[$i] + "_FV2210"

    $rightAddr = $rightCols[$i] + "1"
    $ws.Range($rightAddr).Value = $baseNames[$i] + "_FV2304"
}

# --- 2. Turn the used range into an Excel Table named "Table1". -----------
$dataRange = $ws.Range("A1:U57")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (row 1) in the worksheet view. --------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null

Write-Host "Renamed headers, created table '$($tbl.Name)' and froze header row."
